$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.333.45'
$ws.Range("E2").Value = '  -0.05%  '

$ws.Range("D3").Value = '1.932.50'
$ws.Range("E3").Value = '  +0.02%  '

$ws.Range("D4").Formula = '="1.001"'
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").Formula = '="0.7505"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E5").Value = '  +5.06%  '

$ws.Range("D6").Formula = '="245.05"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E6").Value = '  -2.42%  '

$ws.Range("D7").Formula = '="1.002"'
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").Formula = '="0.3179"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E8").Value = '  -2.61%  '

$ws.Range("D9").Formula = '="27.49"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E9").Value = '  -0.21%  '

$ws.Range("D10").Formula = '="0.06978"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E10").Value = '  -3.05%  '

$ws.Range("D11").Formula = '="0.7791"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E11").Value = '  -2.93%  '

$ws.Range("D12").Formula = '="0.07991"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E12").Value = '  -1.10%  '

$ws.Range("D13").Value = '1.934.53'
$ws.Range("E13").Value = '  +0.18%  '

$ws.Range("D14").Formula = '="5.342"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E14").Value = '  -1.41%  '

$ws.Range("D15").Formula = '="94.26"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E15").Value = '  -0.22%  '

$ws.Range("D16").Formula = '="14.36"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E16").Value = '  -3.87%  '

$ws.Range("D17").Value = '30.349.86'
$ws.Range("E17").Value = '  +0.05%  '

$ws.Range("D18").Formula = '="252.72"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E18").Value = '  +0.00%  '

$ws.Range("D19").Formula = '="0.000007915"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E19").Value = '  -3.30%  '

$ws.Range("D20").Formula = '="5.713"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E20").Value = '  -1.56%  '

$ws.Range("D21").Value = '2.190.14'
$ws.Range("E21").Value = '  +0.14%  '

$ws.Range("D22").Formula = '="1.001"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E22").Value = '  -0.01%  '

$ws.Range("D23").Formula = '="1.001"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E23").Value = '  -0.27%  '

$ws.Range("D24").Formula = '="6.674"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("D25").Formula = '="9.466"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E25").Value = '  -2.34%  '

$ws.Range("D26").Formula = '="166.03"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E26").Value = '  +0.12%  '

$ws.Range("D27").Formula = '="18.94"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E27").Value = '  -1.58%  '

$ws.Range("D28").Formula = '="0.1330"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E28").Value = '  +3.29%  '

$ws.Range("D29").Formula = '="2.216"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E29").Value = '  -5.08%  '

$ws.Range("D30").Formula = '="1.362"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E30").Value = '  -0.21%  '

$ws.Range("D31").Formula = '="1.509"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E31").Value = '  -2.30%  '

$ws.Range("D32").Formula = '="4.370"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E32").Value = '  -1.36%  '

$ws.Range("D33").Formula = '="4.107"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E33").Value = '  -2.18%  '

$ws.Range("D34").Formula = '="0.05154"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E34").Value = '  -1.29%  '

$ws.Range("D35").Formula = '="1.269"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E35").Value = '  +0.00%  '

$ws.Range("D36").Formula = '="0.7449"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E36").Value = '  -0.47%  '

$ws.Range("D37").Formula = '="2.770"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E37").Value = '  +0.13%  '

$ws.Range("D38").Formula = '="0.01945"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E38").Value = '  -1.00%  '

$ws.Range("D39").Formula = '="2.790"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E39").Value = '  -0.51%  '

$ws.Range("D40").Formula = '="77.78"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E40").Value = '  -1.66%  '

$ws.Range("D41").Formula = '="6.393"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E41").Value = '  -1.00%  '

$ws.Range("D42").Formula = '="0.4456"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E42").Value = '  -1.53%  '

$ws.Range("D43").Formula = '="1.960"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E43").Value = '  -3.32%  '

$ws.Range("D44").Formula = '="1.002"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E44").Value = '  +0.02%  '

$ws.Range("E45").Value = '  -1.11%  '

$ws.Range("D46").Formula = '="100.97"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E46").Value = '  -0.94%  '

$ws.Range("D47").Formula = '="9.734"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E47").Value = '  -0.35%  '

$ws.Range("D48").Formula = '="7.455"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E48").Value = '  +0.29%  '

$ws.Range("D49").Formula = '="984.31"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E49").Value = '  +11.36%  '

$ws.Range("D50").Formula = '="37.16"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E50").Value = '  +1.34%  '

$ws.Range("D51").Formula = '="0.06008"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E51").Value = '  -0.84%  '
